$d = $word.ActiveDocument

# Robustly locate the paragraph that ends with "...extract the value from
# there" (the last bullet in the list at the time of this edit) by searching
# for a short, unique snippet of its text and then expanding the found
# range out to the whole enclosing paragraph.
$rng = $d.Content
$found = $rng.Find.Execute("extract the value from there", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not locate the target paragraph (anchor text not found)."
}
$rng.Expand(4) | Out-Null   # wdParagraph: grow the hit out to its full paragraph

# Range covering that paragraph's content, excluding the trailing paragraph mark
$full = $d.Range($rng.Start, $rng.End - 1)

# Replace that paragraph's content (re-asserting its own text, but splitting
# the final run into "...from ther" + "e"), then append the four new bullet
# paragraphs that follow it, all using the same ListParagraph / numId=11
# bullet-list formatting and 12pt (sz/szCs=24) run formatting as the rest of
# the list.
$xml = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="11"/></w:numPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Control components vs Uncontroll components: Ex. Saving the input value into a state property and overwriting the value back to the input. After that we can control what is shown for the user(like only capitalized text).</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> Uncontrol components means that if we want a value we need to reach the DOM and extract the value from ther</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>e</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="11"/></w:numPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>To access the props from a class component we need to refer it with this.props….</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> We can pass a prop from child to parent when passing the prop from parent to child, and then calling the method from the child as this.props.onSubmitParent(this.state.term)</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="11"/></w:numPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Using axios instead of fetch for making request to an API</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="11"/></w:numPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>For list of elements we should add a key</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> to the root returned element</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> so that when the react render our content the performance will increase</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="11"/></w:numPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>React Refs: gives access to a single DOM element; -We create refs in the constructor, assign them to instance variables, then pass to a particular JSX element as props</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>(&lt;img ref={this.imageRef}</w:t></w:r></w:p>'

[void]$full.InsertXML($xml)
